# Insert a new record row at row 232 (pushing existing rows 232..319 down to 233..320)
# and populate the newly inserted row with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(232).Insert()

$ws.Cells.Item(232, 1).Value  = 4
$ws.Cells.Item(232, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(232, 3).Value  = "Los Lagos"
$ws.Cells.Item(232, 4).Value  = 44900
$ws.Cells.Item(232, 5).Value  = 10
$ws.Cells.Item(232, 6).Value  = 100112032
$ws.Cells.Item(232, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(232, 8).Value  = "Sin especificar"
$ws.Cells.Item(232, 9).Value  = "Primera"
$ws.Cells.Item(232, 10).Value = 70
$ws.Cells.Item(232, 11).Value = 15000
$ws.Cells.Item(232, 12).Value = 15000
$ws.Cells.Item(232, 13).Value = 15000
$ws.Cells.Item(232, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(232, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(232, 16).Value = 300
$ws.Cells.Item(232, 17).Value = 50
$ws.Cells.Item(232, 18).Value = "Hortaliza"
